# Update recalculated TPM-based values for Vegfa-Kdr LR-pair sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.628848666666667
$ws.Range("H2").Value = 7.886546000000001
$ws.Range("I2").Value = 0.04622248078033103
$ws.Range("J2").Value = 0.04850184447997802
$ws.Range("M2").Value = 90.43008666666667
$ws.Range("N2").Value = 271.29026
$ws.Range("O2").Value = 0.863466363695901
$ws.Range("P2").Value = 0.8656179140344247
$ws.Range("Q2").Value = 237.7270127602178
$ws.Range("R2").Value = 2139.54311484196
$ws.Range("S2").Value = 0.03991155740039611
$ws.Range("T2").Value = 0.04198406544558065
# Row 3
$ws.Range("G3").Value = 2.628848666666667
$ws.Range("H3").Value = 7.886546000000001
$ws.Range("I3").Value = 0.04622248078033103
$ws.Range("J3").Value = 0.04850184447997802
$ws.Range("M3").Value = 0.06306133333333333
$ws.Range("N3").Value = 0.189184
$ws.Range("O3").Value = 0.000602137432244878
$ws.Range("P3").Value = 0.0006036378137891445
$ws.Range("Q3").Value = 0.1657787020515556
$ws.Range("R3").Value = 1.492008318464
$ws.Range("S3").Value = 0.00002783228588905675
$ws.Range("T3").Value = 0.00002927754736663502
# Row 4
$ws.Range("G4").Value = 2.628848666666667
$ws.Range("H4").Value = 7.886546000000001
$ws.Range("I4").Value = 0.04622248078033103
$ws.Range("J4").Value = 0.04850184447997802
$ws.Range("M4").Value = 9.467965
$ws.Range("N4").Value = 28.403895
$ws.Range("O4").Value = 0.0904043069236993
$ws.Range("P4").Value = 0.09062957269587499
$ws.Range("Q4").Value = 24.88984716629667
$ws.Range("R4").Value = 224.00862449667
$ws.Range("S4").Value = 0.004178711339239838
$ws.Range("T4").Value = 0.004395701440182191
# Row 5
$ws.Range("G5").Value = 2.628848666666667
$ws.Range("H5").Value = 7.886546000000001
$ws.Range("I5").Value = 0.04622248078033103
$ws.Range("J5").Value = 0.04850184447997802
$ws.Range("M5").Value = 0.7809334999999999
$ws.Range("N5").Value = 1.561867
$ws.Range("O5").Value = 0.007456697592460336
$ws.Range("P5").Value = 0.004983518592002547
$ws.Range("Q5").Value = 2.052955990230333
$ws.Range("R5").Value = 12.317735941382
$ws.Range("S5").Value = 0.0003446670611522385
$ws.Range("T5").Value = 0.0002417098437123866
# Row 6
$ws.Range("G6").Value = 2.628848666666667
$ws.Range("H6").Value = 7.886546000000001
$ws.Range("I6").Value = 0.04622248078033103
$ws.Range("J6").Value = 0.04850184447997802
$ws.Range("M6").Value = 3.98709
$ws.Range("N6").Value = 11.96127
$ws.Range("O6").Value = 0.03807049435569441
$ws.Range("P6").Value = 0.03816535686390858
$ws.Range("Q6").Value = 10.48145623038
$ws.Range("R6").Value = 94.33310607342001
$ws.Range("S6").Value = 0.001759712693653786
$ws.Range("T6").Value = 0.001851090203136156
# Row 7
$ws.Range("G7").Value = 26.85202466666667
$ws.Range("H7").Value = 80.556074
$ws.Range("I7").Value = 0.472133375270229
$ws.Range("J7").Value = 0.4954156322762335
$ws.Range("M7").Value = 90.43008666666667
$ws.Range("N7").Value = 271.29026
$ws.Range("O7").Value = 0.863466363695901
$ws.Range("P7").Value = 0.8656179140344247
$ws.Range("Q7").Value = 2428.230917782138
$ws.Range("R7").Value = 21854.07826003924
$ws.Range("S7").Value = 0.4076712887240569
$ws.Range("T7").Value = 0.4288406461909989
# Row 8
$ws.Range("G8").Value = 26.85202466666667
$ws.Range("H8").Value = 80.556074
$ws.Range("I8").Value = 0.472133375270229
$ws.Range("J8").Value = 0.4954156322762335
$ws.Range("M8").Value = 0.06306133333333333
$ws.Range("N8").Value = 0.189184
$ws.Range("O8").Value = 0.000602137432244878
$ws.Range("P8").Value = 0.0006036378137891445
$ws.Range("Q8").Value = 1.693324478179555
$ws.Range("R8").Value = 15.239920303616
$ws.Range("S8").Value = 0.0002842891782623231
$ws.Range("T8").Value = 0.0002990516091841923
# Row 9
$ws.Range("G9").Value = 26.85202466666667
$ws.Range("H9").Value = 80.556074
$ws.Range("I9").Value = 0.472133375270229
$ws.Range("J9").Value = 0.4954156322762335
$ws.Range("M9").Value = 9.467965
$ws.Range("N9").Value = 28.403895
$ws.Range("O9").Value = 0.0904043069236993
$ws.Range("P9").Value = 0.09062957269587499
$ws.Range("Q9").Value = 254.2340297231366
$ws.Range("R9").Value = 2288.10626750823
$ws.Range("S9").Value = 0.04268289056685189
$ws.Range("T9").Value = 0.04489930706005178
# Row 10
$ws.Range("G10").Value = 26.85202466666667
$ws.Range("H10").Value = 80.556074
$ws.Range("I10").Value = 0.472133375270229
$ws.Range("J10").Value = 0.4954156322762335
$ws.Range("M10").Value = 0.7809334999999999
$ws.Range("N10").Value = 1.561867
$ws.Range("O10").Value = 0.007456697592460336
$ws.Range("P10").Value = 0.004983518592002547
$ws.Range("Q10").Value = 20.96964560502633
$ws.Range("R10").Value = 125.817873630158
$ws.Range("S10").Value = 0.003520555802697689
$ws.Range("T10").Value = 0.002468913014217307
# Row 11
$ws.Range("G11").Value = 26.85202466666667
$ws.Range("H11").Value = 80.556074
$ws.Range("I11").Value = 0.472133375270229
$ws.Range("J11").Value = 0.4954156322762335
$ws.Range("M11").Value = 3.98709
$ws.Range("N11").Value = 11.96127
$ws.Range("O11").Value = 0.03807049435569441
$ws.Range("P11").Value = 0.03816535686390858
$ws.Range("Q11").Value = 107.06143902822
$ws.Range("R11").Value = 963.55295125398
$ws.Range("S11").Value = 0.01797435099836021
$ws.Range("T11").Value = 0.01890771440178136
# Row 12
$ws.Range("G12").Value = 11.96574466666667
$ws.Range("H12").Value = 35.897234
$ws.Range("I12").Value = 0.2103911152781009
$ws.Range("J12").Value = 0.2207661073338543
$ws.Range("M12").Value = 90.43008666666667
$ws.Range("N12").Value = 271.29026
$ws.Range("O12").Value = 0.863466363695901
$ws.Range("P12").Value = 0.8656179140344247
$ws.Range("Q12").Value = 1082.063327237871
$ws.Range("R12").Value = 9738.569945140838
$ws.Range("S12").Value = 0.1816656512631069
$ws.Range("T12").Value = 0.1910990973198309
# Row 13
$ws.Range("G13").Value = 11.96574466666667
$ws.Range("H13").Value = 35.897234
$ws.Range("I13").Value = 0.2103911152781009
$ws.Range("J13").Value = 0.2207661073338543
$ws.Range("M13").Value = 0.06306133333333333
$ws.Range("N13").Value = 0.189184
$ws.Range("O13").Value = 0.000602137432244878
$ws.Range("P13").Value = 0.0006036378137891445
$ws.Range("Q13").Value = 0.7545758130062221
$ws.Range("R13").Value = 6.791182317055999
$ws.Range("S13").Value = 0.0001266843659206918
$ws.Range("T13").Value = 0.0001332627703897474
# Row 14
$ws.Range("G14").Value = 11.96574466666667
$ws.Range("H14").Value = 35.897234
$ws.Range("I14").Value = 0.2103911152781009
$ws.Range("J14").Value = 0.2207661073338543
$ws.Range("M14").Value = 9.467965
$ws.Range("N14").Value = 28.403895
$ws.Range("O14").Value = 0.0904043069236993
$ws.Range("P14").Value = 0.09062957269587499
$ws.Range("Q14").Value = 113.2912517029366
$ws.Range("R14").Value = 1019.62126532643
$ws.Range("S14").Value = 0.01902026295962083
$ws.Range("T14").Value = 0.02000793797339889
# Row 15
$ws.Range("G15").Value = 11.96574466666667
$ws.Range("H15").Value = 35.897234
$ws.Range("I15").Value = 0.2103911152781009
$ws.Range("J15").Value = 0.2207661073338543
$ws.Range("M15").Value = 0.7809334999999999
$ws.Range("N15").Value = 1.561867
$ws.Range("O15").Value = 0.007456697592460336
$ws.Range("P15").Value = 0.004983518592002547
$ws.Range("Q15").Value = 9.344450862646331
$ws.Range("R15").Value = 56.06670517587799
$ws.Range("S15").Value = 0.00156882292276926
$ws.Range("T15").Value = 0.001100192000382293
# Row 16
$ws.Range("G16").Value = 11.96574466666667
$ws.Range("H16").Value = 35.897234
$ws.Range("I16").Value = 0.2103911152781009
$ws.Range("J16").Value = 0.2207661073338543
$ws.Range("M16").Value = 3.98709
$ws.Range("N16").Value = 11.96127
$ws.Range("O16").Value = 0.03807049435569441
$ws.Range("P16").Value = 0.03816535686390858
$ws.Range("Q16").Value = 47.70850090302
$ws.Range("R16").Value = 429.37650812718
$ws.Range("S16").Value = 0.008009693766683193
$ws.Range("T16").Value = 0.008425617269852494
# Row 17
$ws.Range("G17").Value = 8.018423
$ws.Range("H17").Value = 16.036846
$ws.Range("I17").Value = 0.1409862072722574
$ws.Range("J17").Value = 0.09862576223372788
$ws.Range("M17").Value = 90.43008666666667
$ws.Range("N17").Value = 271.29026
$ws.Range("O17").Value = 0.863466363695901
$ws.Range("P17").Value = 0.8656179140344247
$ws.Range("Q17").Value = 725.1066868199933
$ws.Range("R17").Value = 4350.64012091996
$ws.Range("S17").Value = 0.1217368477246527
$ws.Range("T17").Value = 0.08537222657481466
# Row 18
$ws.Range("G18").Value = 8.018423
$ws.Range("H18").Value = 16.036846
$ws.Range("I18").Value = 0.1409862072722574
$ws.Range("J18").Value = 0.09862576223372788
$ws.Range("M18").Value = 0.06306133333333333
$ws.Range("N18").Value = 0.189184
$ws.Range("O18").Value = 0.000602137432244878
$ws.Range("P18").Value = 0.0006036378137891445
$ws.Range("Q18").Value = 0.5056524456106667
$ws.Range("R18").Value = 3.033914673664
$ws.Range("S18").Value = 0.00008489307282886124
$ws.Range("T18").Value = 0.00005953423949805547
# Row 19
$ws.Range("G19").Value = 8.018423
$ws.Range("H19").Value = 16.036846
$ws.Range("I19").Value = 0.1409862072722574
$ws.Range("J19").Value = 0.09862576223372788
$ws.Range("M19").Value = 9.467965
$ws.Range("N19").Value = 28.403895
$ws.Range("O19").Value = 0.0904043069236993
$ws.Range("P19").Value = 0.09062957269587499
$ws.Range("Q19").Value = 75.918148319195
$ws.Range("R19").Value = 455.50888991517
$ws.Range("S19").Value = 0.01274576035424945
$ws.Range("T19").Value = 0.008938410688047722
# Row 20
$ws.Range("G20").Value = 8.018423
$ws.Range("H20").Value = 16.036846
$ws.Range("I20").Value = 0.1409862072722574
$ws.Range("J20").Value = 0.09862576223372788
$ws.Range("M20").Value = 0.7809334999999999
$ws.Range("N20").Value = 1.561867
$ws.Range("O20").Value = 0.007456697592460336
$ws.Range("P20").Value = 0.004983518592002547
$ws.Range("Q20").Value = 6.2618551378705
$ws.Range("R20").Value = 25.047420551482
$ws.Range("S20").Value = 0.001051291512337156
$ws.Range("T20").Value = 0.0004915033197422055
# Row 21
$ws.Range("G21").Value = 8.018423
$ws.Range("H21").Value = 16.036846
$ws.Range("I21").Value = 0.1409862072722574
$ws.Range("J21").Value = 0.09862576223372788
$ws.Range("M21").Value = 3.98709
$ws.Range("N21").Value = 11.96127
$ws.Range("O21").Value = 0.03807049435569441
$ws.Range("P21").Value = 0.03816535686390858
$ws.Range("Q21").Value = 31.97017415907
$ws.Range("R21").Value = 191.82104495442
$ws.Range("S21").Value = 0.005367414608189239
$ws.Range("T21").Value = 0.003764087411625222
# Row 22
$ws.Range("G22").Value = 7.408770666666666
$ws.Range("H22").Value = 22.226312
$ws.Range("I22").Value = 0.1302668213990815
$ws.Range("J22").Value = 0.1366906536762062
$ws.Range("M22").Value = 90.43008666666667
$ws.Range("N22").Value = 271.29026
$ws.Range("O22").Value = 0.863466363695901
$ws.Range("P22").Value = 0.8656179140344247
$ws.Range("Q22").Value = 669.9757734801244
$ws.Range("R22").Value = 6029.78196132112
$ws.Range("S22").Value = 0.1124810185836883
$ws.Range("T22").Value = 0.1183218785031996
# Row 23
$ws.Range("G23").Value = 7.408770666666666
$ws.Range("H23").Value = 22.226312
$ws.Range("I23").Value = 0.1302668213990815
$ws.Range("J23").Value = 0.1366906536762062
$ws.Range("M23").Value = 0.06306133333333333
$ws.Range("N23").Value = 0.189184
$ws.Range("O23").Value = 0.000602137432244878
$ws.Range("P23").Value = 0.0006036378137891445
$ws.Range("Q23").Value = 0.4672069566008888
$ws.Range("R23").Value = 4.204862609408
$ws.Range("S23").Value = 0.00007843852934394507
$ws.Range("T23").Value = 0.0000825116473505142
# Row 24
$ws.Range("G24").Value = 7.408770666666666
$ws.Range("H24").Value = 22.226312
$ws.Range("I24").Value = 0.1302668213990815
$ws.Range("J24").Value = 0.1366906536762062
$ws.Range("M24").Value = 9.467965
$ws.Range("N24").Value = 28.403895
$ws.Range("O24").Value = 0.0904043069236993
$ws.Range("P24").Value = 0.09062957269587499
$ws.Range("Q24").Value = 70.14598136502666
$ws.Range("R24").Value = 631.31383228524
$ws.Range("S24").Value = 0.01177668170373729
$ws.Range("T24").Value = 0.0123882155341944
# Row 25
$ws.Range("G25").Value = 7.408770666666666
$ws.Range("H25").Value = 22.226312
$ws.Range("I25").Value = 0.1302668213990815
$ws.Range("J25").Value = 0.1366906536762062
$ws.Range("M25").Value = 0.7809334999999999
$ws.Range("N25").Value = 1.561867
$ws.Range("O25").Value = 0.007456697592460336
$ws.Range("P25").Value = 0.004983518592002547
$ws.Range("Q25").Value = 5.785757207417332
$ws.Range("R25").Value = 34.714543244504
$ws.Range("S25").Value = 0.0009713602935039917
$ws.Range("T25").Value = 0.0006812004139483548
# Row 26
$ws.Range("G26").Value = 7.408770666666666
$ws.Range("H26").Value = 22.226312
$ws.Range("I26").Value = 0.1302668213990815
$ws.Range("J26").Value = 0.1366906536762062
$ws.Range("M26").Value = 3.98709
$ws.Range("N26").Value = 11.96127
$ws.Range("O26").Value = 0.03807049435569441
$ws.Range("P26").Value = 0.03816535686390858
$ws.Range("Q26").Value = 29.53943543736
$ws.Range("R26").Value = 265.85491893624
$ws.Range("S26").Value = 0.004959322288807986
$ws.Range("T26").Value = 0.005216847577513347
